$d = $word.ActiveDocument

function Insert-PlainText($rng, [string]$text) {
    # Appends $text right after the (collapsed) range and advances the range
    # past it, so repeated calls build up a paragraph left-to-right.
    $rng.InsertAfter($text)
    $rng.Collapse(0)
}

function Insert-ParaBreak($rng) {
    $rng.InsertAfter("`r")
    $rng.Collapse(0)
}

# --- Locate the end of the paragraph that currently ends "...on the balance. " ---
$find = $d.Content
$found = $find.Find.Execute("on the balance. ")
if (-not $found) {
    throw "Anchor text not found"
}
$find.Collapse(0)

# --- Continue paragraph 1 with five more sentences (plain Times New Roman) ---
Insert-PlainText $find "This, combined with the first order error calculations made within the LabView software used to interpret the data collected by the balance, would give an uncorrected representation of the flight coefficients at each angle of attack. In order to obtain a full representation of the flight characteristics of the aircraft, it is necessary to include several other corrections for each coefficient obtained during a given experiment."

# --- Start a brand-new paragraph (justified) with the blockage-correction discussion ---
Insert-ParaBreak $find

# Insert ALL the plain text of the paragraph first (so every run inherits the
# paragraph's default "Times New Roman" rPr cleanly), remembering the offsets
# of the spans that need special character formatting, then go back and
# apply that formatting retroactively -- this avoids a quirk where editing a
# sub-range's Font strips rFonts from whatever gets typed immediately after.

Insert-PlainText $find "The first error correction that is usually calculated is that associated with solid blockage of the wing, body, and tail of the aircraft. Each of these aspects of the aircraft disrupt the freestream flow running through the test section, thus interfering with the local flow around the aircraft, and thus, adding to the dynamic pressure that the aircraft experiences. This correction factor is calculated based on the cross-sectional area of the test section, "
# NOTE: must set alignment via the SAME Range object that just performed the
# InsertAfter (a freshly constructed Range over the identical offsets does
# not reliably commit ParagraphFormat changes in this runtime).
$find.ParagraphFormat.Alignment = 3
$spanCEnd         = $find.Start

Insert-PlainText $find "C"
$spanCStart       = $spanCEnd
$spanCEndPos      = $find.Start

Insert-PlainText $find ", the volume associated with the maximum frontal area of the aircraft component being accounted for, and three coefficients (K"
$spanK1Start      = $find.Start
Insert-PlainText $find "1"
$spanK1End        = $find.Start

Insert-PlainText $find ", K"
$spanK2Start      = $find.Start
Insert-PlainText $find "2"
$spanK2End        = $find.Start

Insert-PlainText $find ", "
$spanTauStart     = $find.Start
$tauChar          = [char]0x03C4
Insert-PlainText $find $tauChar
$spanTauEnd       = $find.Start

$spanTau1Start    = $find.Start
Insert-PlainText $find "1"
$spanTau1End      = $find.Start

Insert-PlainText $find ") that are found based on the graphs listed in Figures 3-2 and 3-3 below. "

# --- Now retroactively apply the special character formatting ---
$rC = $d.Range($spanCStart, $spanCEndPos)
$rC.Font.Italic = $true

$rK1 = $d.Range($spanK1Start, $spanK1End)
$rK1.Font.Subscript = $true

$rK2 = $d.Range($spanK2Start, $spanK2End)
$rK2.Font.Subscript = $true

$rTau = $d.Range($spanTauStart, $spanTauEnd)
$rTau.Font.Name = "Calibri"

$rTau1 = $d.Range($spanTau1Start, $spanTau1End)
$rTau1.Font.Subscript = $true

# --- New empty centered paragraph (mirrors the pre-existing one that follows) ---
Insert-ParaBreak $find
$placeholderStart = $find.Start
$find.InsertAfter("X")
$find.ParagraphFormat.Alignment = 1
$find.Collapse(0)
$placeholder = $d.Range($placeholderStart, $placeholderStart + 1)
$placeholder.Delete()
